# Generate Report for Handback
# This script updates the localization-status report after a handback:
#  - Status text moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Handback DateTime stamps are refreshed for zh-cn and de-de
#  - The stale "handback not latest" error message is cleared (report is now in sync)
#  - A couple of report columns are widened/narrowed to fit the new content

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status column for both locales ---
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$ws2.Range("C2").Value = $newStatus
$ws2.Range("K2").Value = "2016-09-07 11:09:10"
$ws2.Range("P2").Value = ""

# --- de-de sheet ---
$ws3.Range("C2").Value = $newStatus
$ws3.Range("K2").Value = "2016-09-07 11:09:20"
$ws3.Range("P2").Value = ""

# --- Column width adjustments to fit the new report content ---
$ws1.Range("E1").ColumnWidth = 29.166666666666668
$ws1.Range("F1").ColumnWidth = 29.166666666666668

$ws2.Range("C1").ColumnWidth = 29.166666666666668
$ws2.Range("P1").ColumnWidth = 12.833333333333334

$ws3.Range("C1").ColumnWidth = 29.166666666666668
$ws3.Range("P1").ColumnWidth = 12.833333333333334
